$wb = $excel.ActiveWorkbook

$edits = @(
    @{ Sheet = "ALC"; Row = 52; Col = "H"; Value = 5244 },
    @{ Sheet = "ALC"; Row = 52; Col = "I"; Value = 5244 },
    @{ Sheet = "ALC"; Row = 52; Col = "J"; Value = 0 },
    @{ Sheet = "ALC"; Row = 52; Col = "K"; Value = 15732 },
    @{ Sheet = "ALC"; Row = 52; Col = "L"; Value = 0 },
    @{ Sheet = "ALC"; Row = 52; Col = "M"; Value = -15572 },
    @{ Sheet = "ALC"; Row = 52; Col = "N"; Value = $null },
    @{ Sheet = "ALC"; Row = 74; Col = "H"; Value = 2894.4546 },
    @{ Sheet = "ALC"; Row = 74; Col = "I"; Value = 1473.1666 },
    @{ Sheet = "ALC"; Row = 74; Col = "K"; Value = 1473.1666 },
    @{ Sheet = "ALC"; Row = 74; Col = "M"; Value = -537.1666 },
    @{ Sheet = "ALC"; Row = 77; Col = "H"; Value = 2894.4546 },
    @{ Sheet = "ALC"; Row = 77; Col = "I"; Value = 1473.1666 },
    @{ Sheet = "ALC"; Row = 77; Col = "K"; Value = 7365.833000000001 },
    @{ Sheet = "ALC"; Row = 77; Col = "M"; Value = -2685.833000000001 },
    @{ Sheet = "ALC"; Row = 100; Col = "H"; Value = 1634.8667 },
    @{ Sheet = "ALC"; Row = 100; Col = "I"; Value = 1277.3 },
    @{ Sheet = "ALC"; Row = 100; Col = "J"; Value = 2350 },
    @{ Sheet = "ALC"; Row = 100; Col = "K"; Value = 1277.3 },
    @{ Sheet = "ALC"; Row = 100; Col = "L"; Value = 2350 },
    @{ Sheet = "ALC"; Row = 100; Col = "M"; Value = -736.3 },
    @{ Sheet = "ALC"; Row = 100; Col = "N"; Value = -3432 },
    @{ Sheet = "ALC"; Row = 116; Col = "H"; Value = 26576084 },
    @{ Sheet = "ALC"; Row = 116; Col = "I"; Value = 25113512 },
    @{ Sheet = "ALC"; Row = 116; Col = "J"; Value = 27794894 },
    @{ Sheet = "ALC"; Row = 116; Col = "K"; Value = 25113512 },
    @{ Sheet = "ALC"; Row = 116; Col = "L"; Value = 27794894 },
    @{ Sheet = "ALC"; Row = 116; Col = "M"; Value = -25110070 },
    @{ Sheet = "ALC"; Row = 116; Col = "N"; Value = -27801778 },
    @{ Sheet = "ALC"; Row = 132; Col = "H"; Value = 3415.3845 },
    @{ Sheet = "ALC"; Row = 132; Col = "I"; Value = 3162.3865 },
    @{ Sheet = "ALC"; Row = 132; Col = "K"; Value = 9487.1595 },
    @{ Sheet = "ALC"; Row = 132; Col = "M"; Value = -6957.1595 },
    @{ Sheet = "ALC"; Row = 137; Col = "H"; Value = 2854.4 },
    @{ Sheet = "ALC"; Row = 137; Col = "I"; Value = 2507.1428 },
    @{ Sheet = "ALC"; Row = 137; Col = "J"; Value = 3664.6667 },
    @{ Sheet = "ALC"; Row = 137; Col = "K"; Value = 7521.428400000001 },
    @{ Sheet = "ALC"; Row = 137; Col = "L"; Value = 10994.0001 },
    @{ Sheet = "ALC"; Row = 137; Col = "M"; Value = -4971.428400000001 },
    @{ Sheet = "ALC"; Row = 137; Col = "N"; Value = -16094.0001 },
    @{ Sheet = "ALC"; Row = 141; Col = "H"; Value = 2423 },
    @{ Sheet = "ALC"; Row = 141; Col = "I"; Value = 1672.5555 },
    @{ Sheet = "ALC"; Row = 141; Col = "K"; Value = 5017.666499999999 },
    @{ Sheet = "ALC"; Row = 141; Col = "M"; Value = 162.3335000000006 },
    @{ Sheet = "ARM"; Row = 33; Col = "H"; Value = 30000 },
    @{ Sheet = "ARM"; Row = 33; Col = "I"; Value = 10000 },
    @{ Sheet = "ARM"; Row = 33; Col = "K"; Value = 10000 },
    @{ Sheet = "ARM"; Row = 33; Col = "M"; Value = -9671 },
    @{ Sheet = "ARM"; Row = 132; Col = "H"; Value = 52633944 },
    @{ Sheet = "ARM"; Row = 132; Col = "I"; Value = 58825876 },
    @{ Sheet = "ARM"; Row = 132; Col = "J"; Value = 2507 },
    @{ Sheet = "ARM"; Row = 132; Col = "K"; Value = 176477628 },
    @{ Sheet = "ARM"; Row = 132; Col = "L"; Value = 7521 },
    @{ Sheet = "ARM"; Row = 132; Col = "M"; Value = -176475098 },
    @{ Sheet = "ARM"; Row = 132; Col = "N"; Value = -12581 },
    @{ Sheet = "BSM"; Row = 22; Col = "H"; Value = 150.25 },
    @{ Sheet = "BSM"; Row = 22; Col = "I"; Value = 182.4 },
    @{ Sheet = "BSM"; Row = 22; Col = "K"; Value = 182.4 },
    @{ Sheet = "BSM"; Row = 22; Col = "M"; Value = -9.400000000000006 },
    @{ Sheet = "CRP"; Row = 21; Col = "H"; Value = 14000 },
    @{ Sheet = "CRP"; Row = 21; Col = "J"; Value = 0 },
    @{ Sheet = "CRP"; Row = 21; Col = "L"; Value = 0 },
    @{ Sheet = "CRP"; Row = 21; Col = "N"; Value = $null },
    @{ Sheet = "CRP"; Row = 31; Col = "H"; Value = 3120.3 },
    @{ Sheet = "CRP"; Row = 31; Col = "I"; Value = 1359.8889 },
    @{ Sheet = "CRP"; Row = 31; Col = "K"; Value = 1359.8889 },
    @{ Sheet = "CRP"; Row = 31; Col = "M"; Value = -1064.8889 },
    @{ Sheet = "CRP"; Row = 34; Col = "H"; Value = 3120.3 },
    @{ Sheet = "CRP"; Row = 34; Col = "I"; Value = 1359.8889 },
    @{ Sheet = "CRP"; Row = 34; Col = "K"; Value = 1359.8889 },
    @{ Sheet = "CRP"; Row = 34; Col = "M"; Value = -1157.8889 },
    @{ Sheet = "CRP"; Row = 50; Col = "H"; Value = 70987 },
    @{ Sheet = "CRP"; Row = 50; Col = "I"; Value = 0 },
    @{ Sheet = "CRP"; Row = 50; Col = "J"; Value = 70987 },
    @{ Sheet = "CRP"; Row = 50; Col = "K"; Value = 0 },
    @{ Sheet = "CRP"; Row = 50; Col = "L"; Value = 70987 },
    @{ Sheet = "CRP"; Row = 50; Col = "M"; Value = $null },
    @{ Sheet = "CRP"; Row = 50; Col = "N"; Value = -72237 },
    @{ Sheet = "CRP"; Row = 51; Col = "H"; Value = 74699 },
    @{ Sheet = "CRP"; Row = 51; Col = "J"; Value = 80097.5 },
    @{ Sheet = "CRP"; Row = 51; Col = "L"; Value = 80097.5 },
    @{ Sheet = "CRP"; Row = 51; Col = "N"; Value = -81569.5 },
    @{ Sheet = "CRP"; Row = 59; Col = "H"; Value = 59999.5 },
    @{ Sheet = "CRP"; Row = 59; Col = "J"; Value = 50000 },
    @{ Sheet = "CRP"; Row = 59; Col = "L"; Value = 50000 },
    @{ Sheet = "CRP"; Row = 59; Col = "N"; Value = -52290 },
    @{ Sheet = "CRP"; Row = 60; Col = "H"; Value = 22723.309 },
    @{ Sheet = "CRP"; Row = 60; Col = "I"; Value = 8703 },
    @{ Sheet = "CRP"; Row = 60; Col = "J"; Value = 25272.455 },
    @{ Sheet = "CRP"; Row = 60; Col = "K"; Value = 8703 },
    @{ Sheet = "CRP"; Row = 60; Col = "L"; Value = 25272.455 },
    @{ Sheet = "CRP"; Row = 60; Col = "M"; Value = -8192 },
    @{ Sheet = "CRP"; Row = 60; Col = "N"; Value = -26294.455 },
    @{ Sheet = "CRP"; Row = 61; Col = "H"; Value = 74699 },
    @{ Sheet = "CRP"; Row = 61; Col = "J"; Value = 80097.5 },
    @{ Sheet = "CRP"; Row = 61; Col = "L"; Value = 80097.5 },
    @{ Sheet = "CRP"; Row = 61; Col = "N"; Value = -80793.5 },
    @{ Sheet = "CRP"; Row = 74; Col = "H"; Value = 35553.375 },
    @{ Sheet = "CRP"; Row = 74; Col = "J"; Value = 35553.375 },
    @{ Sheet = "CRP"; Row = 74; Col = "L"; Value = 35553.375 },
    @{ Sheet = "CRP"; Row = 74; Col = "N"; Value = -37301.375 },
    @{ Sheet = "CRP"; Row = 77; Col = "H"; Value = 35553.375 },
    @{ Sheet = "CRP"; Row = 77; Col = "J"; Value = 35553.375 },
    @{ Sheet = "CRP"; Row = 77; Col = "L"; Value = 106660.125 },
    @{ Sheet = "CRP"; Row = 77; Col = "N"; Value = -115396.125 },
    @{ Sheet = "CRP"; Row = 94; Col = "H"; Value = 1208.5 },
    @{ Sheet = "CRP"; Row = 94; Col = "I"; Value = 1079.2 },
    @{ Sheet = "CRP"; Row = 94; Col = "J"; Value = 1337.8 },
    @{ Sheet = "CRP"; Row = 94; Col = "K"; Value = 1079.2 },
    @{ Sheet = "CRP"; Row = 94; Col = "L"; Value = 1337.8 },
    @{ Sheet = "CRP"; Row = 94; Col = "M"; Value = -628.2 },
    @{ Sheet = "CRP"; Row = 94; Col = "N"; Value = -2239.8 },
    @{ Sheet = "CUL"; Row = 23; Col = "H"; Value = 100 },
    @{ Sheet = "CUL"; Row = 23; Col = "J"; Value = 100 },
    @{ Sheet = "CUL"; Row = 23; Col = "L"; Value = 300 },
    @{ Sheet = "CUL"; Row = 23; Col = "N"; Value = -770 },
    @{ Sheet = "CUL"; Row = 129; Col = "H"; Value = 1871.7646 },
    @{ Sheet = "CUL"; Row = 129; Col = "I"; Value = 798.5 },
    @{ Sheet = "CUL"; Row = 129; Col = "J"; Value = 2825.7778 },
    @{ Sheet = "CUL"; Row = 129; Col = "K"; Value = 2395.5 },
    @{ Sheet = "CUL"; Row = 129; Col = "L"; Value = 8477.3334 },
    @{ Sheet = "CUL"; Row = 129; Col = "M"; Value = 2604.5 },
    @{ Sheet = "CUL"; Row = 129; Col = "N"; Value = -18477.3334 },
    @{ Sheet = "CUL"; Row = 136; Col = "H"; Value = 1111.5 },
    @{ Sheet = "CUL"; Row = 136; Col = "I"; Value = 1111.5 },
    @{ Sheet = "CUL"; Row = 136; Col = "J"; Value = 0 },
    @{ Sheet = "CUL"; Row = 136; Col = "K"; Value = 3334.5 },
    @{ Sheet = "CUL"; Row = 136; Col = "L"; Value = 0 },
    @{ Sheet = "CUL"; Row = 136; Col = "M"; Value = 1765.5 },
    @{ Sheet = "CUL"; Row = 136; Col = "N"; Value = $null },
    @{ Sheet = "CUL"; Row = 139; Col = "H"; Value = 53344790 },
    @{ Sheet = "CUL"; Row = 139; Col = "I"; Value = 66680064 },
    @{ Sheet = "CUL"; Row = 139; Col = "K"; Value = 200040192 },
    @{ Sheet = "CUL"; Row = 139; Col = "M"; Value = -200035052 },
    @{ Sheet = "LTW"; Row = 40; Col = "H"; Value = 2240 },
    @{ Sheet = "LTW"; Row = 40; Col = "I"; Value = 1925 },
    @{ Sheet = "LTW"; Row = 40; Col = "J"; Value = 2450 },
    @{ Sheet = "LTW"; Row = 40; Col = "K"; Value = 1925 },
    @{ Sheet = "LTW"; Row = 40; Col = "L"; Value = 2450 },
    @{ Sheet = "LTW"; Row = 40; Col = "M"; Value = -1789 },
    @{ Sheet = "LTW"; Row = 40; Col = "N"; Value = -2722 },
    @{ Sheet = "LTW"; Row = 46; Col = "H"; Value = 1206.25 },
    @{ Sheet = "LTW"; Row = 46; Col = "I"; Value = 749.5 },
    @{ Sheet = "LTW"; Row = 46; Col = "J"; Value = 1271.5 },
    @{ Sheet = "LTW"; Row = 46; Col = "K"; Value = 749.5 },
    @{ Sheet = "LTW"; Row = 46; Col = "L"; Value = 1271.5 },
    @{ Sheet = "LTW"; Row = 46; Col = "M"; Value = -561.5 },
    @{ Sheet = "LTW"; Row = 46; Col = "N"; Value = -1647.5 },
    @{ Sheet = "LTW"; Row = 61; Col = "H"; Value = 14310.223 },
    @{ Sheet = "LTW"; Row = 61; Col = "I"; Value = 9965.5 },
    @{ Sheet = "LTW"; Row = 61; Col = "K"; Value = 9965.5 },
    @{ Sheet = "LTW"; Row = 61; Col = "M"; Value = -9763.5 },
    @{ Sheet = "LTW"; Row = 113; Col = "H"; Value = 14310.223 },
    @{ Sheet = "LTW"; Row = 113; Col = "I"; Value = 9965.5 },
    @{ Sheet = "LTW"; Row = 113; Col = "K"; Value = 9965.5 },
    @{ Sheet = "LTW"; Row = 113; Col = "M"; Value = -7795.5 },
    @{ Sheet = "LTW"; Row = 122; Col = "H"; Value = 3490.9443 },
    @{ Sheet = "LTW"; Row = 122; Col = "I"; Value = 2774.0715 },
    @{ Sheet = "LTW"; Row = 122; Col = "K"; Value = 8322.2145 },
    @{ Sheet = "LTW"; Row = 122; Col = "M"; Value = -5872.2145 },
    @{ Sheet = "WVR"; Row = 43; Col = "H"; Value = 19014.572 },
    @{ Sheet = "WVR"; Row = 43; Col = "I"; Value = 19008.5 },
    @{ Sheet = "WVR"; Row = 43; Col = "K"; Value = 19008.5 },
    @{ Sheet = "WVR"; Row = 43; Col = "M"; Value = -18859.5 },
    @{ Sheet = "WVR"; Row = 81; Col = "H"; Value = 9531873 },
    @{ Sheet = "WVR"; Row = 81; Col = "I"; Value = 6523 },
    @{ Sheet = "WVR"; Row = 81; Col = "K"; Value = 13046 },
    @{ Sheet = "WVR"; Row = 81; Col = "M"; Value = -11985 },
    @{ Sheet = "WVR"; Row = 84; Col = "H"; Value = 9531873 },
    @{ Sheet = "WVR"; Row = 84; Col = "I"; Value = 6523 },
    @{ Sheet = "WVR"; Row = 84; Col = "K"; Value = 65230 },
    @{ Sheet = "WVR"; Row = 84; Col = "M"; Value = -59926 }
)

foreach ($e in $edits) {
    $ws = $wb.Worksheets.Item($e.Sheet)
    $addr = "$($e.Col)$($e.Row)"
    if ($null -eq $e.Value) {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value = $e.Value
    }
}

Write-Output ("Applied " + $edits.Count + " cell edits.")
